# Auto-generated Excel COM-interop script applying the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)

# Insert a brand-new row at position 42 (shifts rows 42-44 down to 43-45).
$ws1.Rows.Item(42).Insert()

# Copy the bold/bordered number-column style from the row above onto the new A42 cell.
$ws1.Cells.Item(41, 1).Copy()
$ws1.Cells.Item(42, 1).PasteSpecial(-4122)

# Fill in the new row 42 content (keeps the same running index as old row 42, i.e. 41).
$ws1.Cells.Item(42, 1).Value = 41
$ws1.Cells.Item(42, 2).NumberFormat = "@"
$ws1.Cells.Item(42, 2).Value = "2024-09-15"
$ws1.Cells.Item(42, 3).Value = "杭州·西溪银泰 布谷布Goods二次元吃谷嘉年华 免票"
$ws1.Cells.Item(42, 4).Value = "双龙街588号 西溪银泰城"
$ws1.Cells.Item(42, 5).Value = "2024.09.15 10:00-09.17 20:00"
$ws1.Cells.Item(42, 6).Value = 0
$ws1.Cells.Item(42, 7).Value = 30
$ws1.Cells.Item(42, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89782"
$ws1.Cells.Item(42, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/iWnJHkey1721737159663.png"

# The event that used to be row 42 is now row 43; its "want to go" count also rose.
$ws1.Cells.Item(43, 6).Value = 853
# The event that used to be row 43 is now row 44; its "want to go" count also rose.
$ws1.Cells.Item(44, 6).Value = 78
# The event that used to be row 44 is now row 45; its "want to go" count also rose.
$ws1.Cells.Item(45, 6).Value = 397

# Plain refresh-count ("想去人数") bumps for the rest of the sheet.
$ws1.Cells.Item(3, 6).Value = 1760
$ws1.Cells.Item(4, 6).Value = 49
$ws1.Cells.Item(6, 6).Value = 2146
$ws1.Cells.Item(7, 6).Value = 1370
$ws1.Cells.Item(8, 6).Value = 2079
$ws1.Cells.Item(9, 6).Value = 967
$ws1.Cells.Item(11, 6).Value = 2404
$ws1.Cells.Item(12, 6).Value = 659
$ws1.Cells.Item(13, 6).Value = 837
$ws1.Cells.Item(14, 6).Value = 3925
$ws1.Cells.Item(15, 6).Value = 316
$ws1.Cells.Item(16, 6).Value = 366
$ws1.Cells.Item(17, 6).Value = 3010
$ws1.Cells.Item(18, 6).Value = 806
$ws1.Cells.Item(20, 6).Value = 1346
$ws1.Cells.Item(21, 6).Value = 117
$ws1.Cells.Item(22, 6).Value = 2050
$ws1.Cells.Item(23, 6).Value = 1178
$ws1.Cells.Item(24, 6).Value = 1882
$ws1.Cells.Item(25, 6).Value = 384
$ws1.Cells.Item(26, 6).Value = 207
$ws1.Cells.Item(28, 6).Value = 8397
$ws1.Cells.Item(29, 6).Value = 5669
$ws1.Cells.Item(32, 6).Value = 751
$ws1.Cells.Item(33, 6).Value = 767
$ws1.Cells.Item(34, 6).Value = 3477
$ws1.Cells.Item(37, 6).Value = 389
$ws1.Cells.Item(38, 6).Value = 36
$ws1.Cells.Item(40, 6).Value = 160
$ws1.Cells.Item(41, 6).Value = 4622

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(15, 6).Value = 107
$ws2.Cells.Item(18, 6).Value = 176
$ws2.Cells.Item(26, 6).Value = 28
$ws2.Cells.Item(3, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/TviSO9CG1721810981388.png"

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 8283
$ws3.Cells.Item(3, 6).Value = 363
$ws3.Cells.Item(4, 6).Value = 1283

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 363
$ws4.Cells.Item(4, 6).Value = 1283
$ws4.Cells.Item(6, 6).Value = 1760
$ws4.Cells.Item(7, 6).Value = 49
$ws4.Cells.Item(9, 6).Value = 1370
$ws4.Cells.Item(10, 6).Value = 2079
$ws4.Cells.Item(11, 6).Value = 967
$ws4.Cells.Item(15, 6).Value = 3925
$ws4.Cells.Item(16, 6).Value = 366
$ws4.Cells.Item(17, 6).Value = 3010
$ws4.Cells.Item(18, 6).Value = 806
$ws4.Cells.Item(21, 6).Value = 2050
$ws4.Cells.Item(27, 6).Value = 1882
$ws4.Cells.Item(28, 6).Value = 107
$ws4.Cells.Item(29, 6).Value = 207
$ws4.Cells.Item(31, 6).Value = 8397
$ws4.Cells.Item(32, 6).Value = 5669
$ws4.Cells.Item(36, 6).Value = 751
$ws4.Cells.Item(37, 6).Value = 767
$ws4.Cells.Item(40, 6).Value = 389
$ws4.Cells.Item(42, 6).Value = 160
$ws4.Cells.Item(43, 6).Value = 4622
$ws4.Cells.Item(44, 6).Value = 853
$ws4.Cells.Item(45, 6).Value = 397
$ws4.Cells.Item(48, 6).Value = 28
$ws4.Cells.Item(5, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/TviSO9CG1721810981388.png"

